# PandaSat CONOPS.pptx - text tweaks on the flowchart decision labels
# 1) "Downlink Requested?" -> "Downlink requested?"                         (TextBox 122)
# 2) "Payload Operation Scheduled?" -> "Payload operation scheduled now?"    (TextBox 143)
# 3) "Battery Voltage above threshold?" -> "Battery voltage above threshold?"(TextBox 90, inside Group 1)
# 4) "Battery Voltage above threshold?" -> split into two runs:
#       "Battery voltage " + "above threshold?"                             (TextBox 292, inside Group 290)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Downlink Requested? -> Downlink requested?
$downlink = $s.Shapes.Item("TextBox 122")
$downlink.TextFrame.TextRange.Text = "Downlink requested?"

# 2) Payload Operation Scheduled? -> Payload operation scheduled now?
$payload = $s.Shapes.Item("TextBox 143")
$payload.TextFrame.TextRange.Text = "Payload operation scheduled now?"

# 3) Battery Voltage above threshold? -> Battery voltage above threshold? (single run, case-only)
$group1 = $s.Shapes.Item("Group 1")
$batt1 = $group1.GroupItems.Item("TextBox 90")
$batt1.TextFrame.TextRange.Text = "Battery voltage above threshold?"

# 4) Battery Voltage above threshold? -> "Battery voltage " + "above threshold?" (two runs)
$group290 = $s.Shapes.Item("Group 290")
$batt2 = $group290.GroupItems.Item("TextBox 292")
$batt2Range = $batt2.TextFrame.TextRange
# Re-case just the leading "Battery Voltage " portion (16 chars incl. trailing space);
# this naturally splits the paragraph into two runs, matching the source edit.
$prefix = $batt2Range.Characters(1, 16)
$prefix.Text = "Battery voltage "
